# Daily attendance processing - 2026-01-04 11:30:01
# Reorders the "Recorded By" email lists (column G) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = "Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    "G3"  = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    "G4"  = "asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    "G5"  = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
    "G6"  = "manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
    "G7"  = "Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
    "G9"  = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
    "G12" = "Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
    "G13" = "yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
    "G15" = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
    "G17" = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
    "G19" = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
    "G20" = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    "G24" = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
    "G25" = "menna-allah.gamil@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg"
    "G27" = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
    "G28" = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
    "G30" = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
